$wb = $excel.ActiveWorkbook

# --- Male_25m: swap rows 8 and 9 (duplicate/mismatched name fix) ---
$ws1 = $wb.Worksheets.Item("Male_25m")
$ws1.Range("A8").Value = "Christian Tronvoll"
$ws1.Range("B8").Value = "2.14,69"
$ws1.Range("D8").Value = "18.01.2014"
$ws1.Range("E8").Value = "Trondheim"

$ws1.Range("A9").Value = "Manith Randula Attanapola"
$ws1.Range("B9").Value = "2.14,67"
$ws1.Range("D9").Value = "29.11.2014"
$ws1.Range("E9").Value = "Stjørdal"

# --- Male_50m: swap rows 5 and 6 (duplicate/mismatched name fix) ---
$ws2 = $wb.Worksheets.Item("Male_50m")
$ws2.Range("A5").Value = "Gabriel Rognes Steen"
$ws2.Range("B5").Value = "2.18,49"
$ws2.Range("D5").Value = "22.05.2021"
$ws2.Range("E5").Value = "Funchal"

$ws2.Range("A6").Value = "Tudor Ignat"
$ws2.Range("B6").Value = "2.18,54"
# "08.08.2023" is ambiguous (day <= 12) and gets auto-parsed as a date by
# the engine's literal-entry heuristics; force literal text via a leading
# apostrophe and then strip the resulting quote-prefix formatting so the
# cell ends up as a plain string value/style, matching the source file.
$ws2.Range("D6").Value = "'08.08.2023"
$ws2.Range("D6").ClearFormats()
$ws2.Range("E6").Value = "Kyushu"
